$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 44523
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B20").Value = 2436
$ws.Range("C20").Value = 0.01
$ws.Range("D20").Value = 40
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = ""
